$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Fix typo: MonteMorelos -> Montemorelos
$ws.Range("B768").Value = "Montemorelos"

# Title-case Spanish connector words (de/del/el/y/la/las/los) in municipality/state names
$ws.Range("B7").Value = 'Pabellón De Arteaga'
$ws.Range("B8").Value = 'Rincón De Romos'
$ws.Range("B36").Value = 'Amatenango De La Frontera'
$ws.Range("B46").Value = 'Chiapa De Corzo'
$ws.Range("B51").Value = 'Comitán De Domínguez'
$ws.Range("B76").Value = 'Marqués De Comillas'
$ws.Range("B77").Value = 'Mazapa De Madero'
$ws.Range("B84").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B93").Value = 'Salto De Agua'
$ws.Range("B94").Value = 'San Cristóbal De Las Casas'
$ws.Range("B129").Value = 'Guadalupe Y Calvo'
$ws.Range("B132").Value = 'Hidalgo Del Parral'
$ws.Range("B140").Value = 'San Francisco Del Oro'
$ws.Range("B143").Value = 'Valle De Zaragoza'
$ws.Range("B160").Value = 'San Juan De Sabinas'
$ws.Range("B174").Value = 'Villa De Álvarez'
$ws.Range("A176").Value = 'Ciudad De México'
$ws.Range("B180").Value = 'Cuajimalpa De Morelos'
$ws.Range("B194").Value = 'Coneto De Comonfort'
$ws.Range("B206").Value = 'Nombre De Dios'
$ws.Range("B210").Value = 'Pánuco De Coronado'
$ws.Range("B216").Value = 'San Juan De Guadalupe'
$ws.Range("B217").Value = 'San Juan Del Río'
$ws.Range("B218").Value = 'San Luis Del Cordero'
$ws.Range("A227").Value = 'Estado De México'
$ws.Range("B227").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B230").Value = 'Almoloya De Alquisiras'
$ws.Range("B231").Value = 'Almoloya De Juárez'
$ws.Range("B232").Value = 'Almoloya Del Río'
$ws.Range("B238").Value = 'Atizapán De Zaragoza'
$ws.Range("B242").Value = 'Chapa De Mota'
$ws.Range("B246").Value = 'Coacalco De Berriozábal'
$ws.Range("B252").Value = 'Ecatepec De Morelos'
$ws.Range("B258").Value = 'Ixtapan De La Sal'
$ws.Range("B259").Value = 'Ixtapan Del Oro'
$ws.Range("B273").Value = 'Naucalpan De Juárez'
$ws.Range("B281").Value = 'San Felipe Del Progreso'
$ws.Range("B283").Value = 'San Simón De Guerrero'
$ws.Range("B285").Value = 'Soyaniquilpan De Juárez'
$ws.Range("B294").Value = 'Tenango Del Valle'
$ws.Range("B306").Value = 'Tlalnepantla De Baz'
$ws.Range("B311").Value = 'Valle De Bravo'
$ws.Range("B312").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B313").Value = 'Villa De Allende'
$ws.Range("B326").Value = 'Apaseo El Alto'
$ws.Range("B327").Value = 'Apaseo El Grande'
$ws.Range("B335").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B339").Value = 'Jaral Del Progreso'
$ws.Range("B347").Value = 'Purísima Del Rincón'
$ws.Range("B351").Value = 'San Diego De La Unión'
$ws.Range("B353").Value = 'San Francisco Del Rincón'
$ws.Range("B355").Value = 'San Luis De La Paz'
$ws.Range("B357").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B358").Value = 'Silao De La Victoria'
$ws.Range("B363").Value = 'Valle De Santiago'
$ws.Range("B369").Value = 'Acapulco De Juárez'
$ws.Range("B371").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B372").Value = 'Alcozauca De Guerrero'
$ws.Range("B376").Value = 'Atenango Del Río'
$ws.Range("B377").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B379").Value = 'Atoyac De Álvarez'
$ws.Range("B380").Value = 'Ayutla De Los Libres'
$ws.Range("B383").Value = 'Buenavista De Cuéllar'
$ws.Range("B384").Value = 'Chilapa De Álvarez'
$ws.Range("B385").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B386").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B391").Value = 'Coyuca De Benítez'
$ws.Range("B392").Value = 'Coyuca De Catalán'
$ws.Range("B396").Value = 'Cuetzala Del Progreso'
$ws.Range("B397").Value = 'Cutzamala De Pinzón'
$ws.Range("B403").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B404").Value = 'Iguala De La Independencia'
$ws.Range("B406").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B407").Value = 'Zihuatanejo De Azueta'
$ws.Range("B409").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B412").Value = 'Mártir De Cuilapan'
$ws.Range("B425").Value = 'Taxco De Alarcón'
$ws.Range("B427").Value = 'Técpan De Galeana'
$ws.Range("B429").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B431").Value = 'Tixtla De Guerrero'
$ws.Range("B435").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B436").Value = 'Tlapa De Comonfort'
$ws.Range("B447").Value = 'Agua Blanca De Iturbide'
$ws.Range("B453").Value = 'Atotonilco De Tula'
$ws.Range("B454").Value = 'Atotonilco El Grande'
$ws.Range("B460").Value = 'Cuautepec De Hinojosa'
$ws.Range("B463").Value = 'Huasca De Ocampo'
$ws.Range("B466").Value = 'Huejutla De Reyes'
$ws.Range("B469").Value = 'Jacala De Ledezma'
$ws.Range("B473").Value = 'Mineral Del Chico'
$ws.Range("B474").Value = 'Mineral Del Monte'
$ws.Range("B475").Value = 'Mixquiahuala De Juárez'
$ws.Range("B476").Value = 'Molango De Escamilla'
$ws.Range("B478").Value = 'Nopala De Villagrán'
$ws.Range("B479").Value = 'Omitlán De Juárez'
$ws.Range("B480").Value = 'Pachuca De Soto'
$ws.Range("B483").Value = 'Progreso De Obregón'
$ws.Range("B489").Value = 'Santiago De Anaya'
$ws.Range("B490").Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range("B494").Value = 'Tenango De Doria'
$ws.Range("B496").Value = 'Tepehuacán De Guerrero'
$ws.Range("B497").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B500").Value = 'Tezontepec De Aldama'
$ws.Range("B509").Value = 'Tula De Allende'
$ws.Range("B510").Value = 'Tulancingo De Bravo'
$ws.Range("B511").Value = 'Villa De Tezontepec'
$ws.Range("B515").Value = 'Zacualtipán De Ángeles'
$ws.Range("B516").Value = 'Zapotlán De Juárez'
$ws.Range("B520").Value = 'Acatlán De Juárez'
$ws.Range("B521").Value = 'Ahualulco De Mercado'
$ws.Range("B525").Value = 'Atotonilco El Alto'
$ws.Range("B527").Value = 'Autlán De Navarro'
$ws.Range("B542").Value = 'Encarnación De Díaz'
$ws.Range("B546").Value = 'Huejuquilla El Alto'
$ws.Range("B547").Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range("B548").Value = 'Ixtlahuacán Del Río'
$ws.Range("B554").Value = 'Lagos De Moreno'
$ws.Range("B560").Value = 'Ojuelos De Jalisco'
$ws.Range("B565").Value = 'San Diego De Alejandría'
$ws.Range("B567").Value = 'San Juan De Los Lagos'
$ws.Range("B569").Value = 'San Martín De Bolaños'
$ws.Range("B571").Value = 'San Sebastián Del Oeste'
$ws.Range("B574").Value = 'Talpa De Allende'
$ws.Range("B575").Value = 'Tamazula De Gordiano'
$ws.Range("B577").Value = 'Techaluta De Montenegro'
$ws.Range("B581").Value = 'Teocuitatlán De Corona'
$ws.Range("B582").Value = 'Tepatitlán De Morelos'
$ws.Range("B584").Value = 'Tizapán El Alto'
$ws.Range("B585").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B593").Value = 'Unión De San Antonio'
$ws.Range("B594").Value = 'Unión De Tula'
$ws.Range("B595").Value = 'Valle De Juárez'
$ws.Range("B600").Value = 'Yahualica De González Gallo'
$ws.Range("B601").Value = 'Zacoalco De Torres'
$ws.Range("B604").Value = 'Zapotitlán De Vadillo'
$ws.Range("B605").Value = 'Zapotlán Del Rey'
$ws.Range("B606").Value = 'Zapotlán El Grande'
$ws.Range("B690").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B714").Value = 'Coatlán Del Río'
$ws.Range("B721").Value = 'Jonacatepec De Leandro Valle'
$ws.Range("B725").Value = 'Puente De Ixtla'
$ws.Range("B731").Value = 'Tlaltizapán De Zapata'
$ws.Range("B741").Value = 'Amatlán De Cañas'
$ws.Range("B744").Value = 'Ixtlán Del Río'
$ws.Range("B751").Value = 'Santa María Del Oro'
$ws.Range("B767").Value = 'Mier Y Noriega'
$ws.Range("B770").Value = 'San Nicolás De Los Garza'
$ws.Range("B774").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B778").Value = 'Ayoquezco De Aldama'
$ws.Range("B783").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B784").Value = 'Chiquihuitlán De Benito Juárez'
$ws.Range("B787").Value = 'Coicoyán De Las Flores'
$ws.Range("B788").Value = 'Constancia Del Rosario'
$ws.Range("B791").Value = 'Guevea De Humboldt'
$ws.Range("B792").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B793").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B794").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B795").Value = 'Huautla De Jiménez'
$ws.Range("B797").Value = 'Ixtlán De Juárez'
$ws.Range("B798").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B806").Value = 'Mártires De Tacubaya'
$ws.Range("B809").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B811").Value = 'Nejapa De Madero'
$ws.Range("B812").Value = 'Oaxaca De Juárez'
$ws.Range("B813").Value = 'Ocotlán De Morelos'
$ws.Range("B814").Value = 'Pinotepa De Don Luis'
$ws.Range("B816").Value = 'Putla Villa De Guerrero'
$ws.Range("B817").Value = 'Reforma De Pineda'
$ws.Range("B834").Value = 'San Antonino El Alto'
$ws.Range("B841").Value = 'San Baltazar Yatzachi El Bajo'
$ws.Range("B851").Value = 'San Felipe Jalapa De Díaz'
$ws.Range("B856").Value = 'San Francisco Del Mar'
$ws.Range("B874").Value = 'San José Del Progreso'
$ws.Range("B882").Value = 'San Juan Bautista Lo De Soto'
$ws.Range("B913").Value = 'San Martín De Los Cansecos'
$ws.Range("B928").Value = 'San Miguel Del Puerto'
$ws.Range("B929").Value = 'San Miguel El Grande'
$ws.Range("B948").Value = 'San Pablo Villa De Mitla'
$ws.Range("B955").Value = 'San Pedro El Alto'
$ws.Range("B990").Value = 'Santa Ana Del Valle'
$ws.Range("B1013").Value = 'Santa Inés De Zaragoza'
$ws.Range("B1014").Value = 'Santa Inés Del Monte'
$ws.Range("B1016").Value = 'Santa Lucía Del Camino'
$ws.Range("B1027").Value = 'Santa María Del Tule'
$ws.Range("B1035").Value = 'Santa María Jalapa Del Marqués'
$ws.Range("B1089").Value = 'Santo Domingo De Morelos'
$ws.Range("B1106").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B1107").Value = 'Tataltepec De Valdés'
$ws.Range("B1108").Value = 'Teococuilco De Marcos Pérez'
$ws.Range("B1109").Value = 'Teotitlán De Flores Magón'
$ws.Range("B1110").Value = 'Teotitlán Del Valle'
$ws.Range("B1112").Value = 'Tezoatlán De Segura Y Luna'
$ws.Range("B1113").Value = 'Tlacolula De Matamoros'
$ws.Range("B1117").Value = 'Villa De Chilapa De Díaz'
$ws.Range("B1118").Value = 'Villa De Etla'
$ws.Range("B1119").Value = 'Villa De Tamazulápam Del Progreso'
$ws.Range("B1120").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B1121").Value = 'Villa De Zaachila'
$ws.Range("B1123").Value = 'Villa Sola De Vega'
$ws.Range("B1124").Value = 'Villa Talea De Castro'
$ws.Range("B1126").Value = 'Zapotitlán Del Río'
$ws.Range("B1129").Value = 'Zimatlán De Álvarez'
$ws.Range("B1152").Value = 'Chalchicomula De Sesma'
$ws.Range("B1161").Value = 'Chila De La Sal'
$ws.Range("B1178").Value = 'Huehuetlán El Chico'
$ws.Range("B1181").Value = 'Huitzilan De Serdán'
$ws.Range("B1183").Value = 'Ixcamilpa De Guerrero'
$ws.Range("B1186").Value = 'Izúcar De Matamoros'
$ws.Range("B1194").Value = 'Los Reyes De Juárez'
$ws.Range("B1203").Value = 'Palmar De Bravo'
$ws.Range("B1222").Value = 'San Salvador El Seco'
$ws.Range("B1223").Value = 'San Salvador El Verde'
$ws.Range("B1229").Value = 'Tecali De Herrera'
$ws.Range("B1235").Value = 'Tepanco De López'
$ws.Range("B1236").Value = 'Tepatlaxco De Hidalgo'
$ws.Range("B1242").Value = 'Tepexi De Rodríguez'
$ws.Range("B1244").Value = 'Tetela De Ocampo'
$ws.Range("B1245").Value = 'Teteles De Avila Castillo'
$ws.Range("B1249").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B1261").Value = 'Totoltepec De Guerrero'
$ws.Range("B1265").Value = 'Xayacatlán De Bravo'
$ws.Range("B1279").Value = 'Amealco De Bonfil'
$ws.Range("B1281").Value = 'Cadereyta De Montes'
$ws.Range("B1287").Value = 'Jalpan De Serra'
$ws.Range("B1288").Value = 'Landa De Matamoros'
$ws.Range("B1291").Value = 'Pinal De Amoles'
$ws.Range("B1294").Value = 'San Juan Del Río'
$ws.Range("B1307").Value = 'Axtla De Terrazas'
$ws.Range("B1311").Value = 'Ciudad Del Maíz'
$ws.Range("B1321").Value = 'Mexquitic De Carmona'
$ws.Range("B1327").Value = 'San Ciro De Acosta'
$ws.Range("B1332").Value = 'Santa María Del Río'
$ws.Range("B1334").Value = 'Soledad De Graciano Sánchez'
$ws.Range("B1341").Value = 'Tanquián De Escobedo'
$ws.Range("B1345").Value = 'Villa De Arista'
$ws.Range("B1346").Value = 'Villa De Arriaga'
$ws.Range("B1347").Value = 'Villa De Guadalupe'
$ws.Range("B1348").Value = 'Villa De Ramos'
$ws.Range("B1349").Value = 'Villa De Reyes'
$ws.Range("B1391").Value = 'Jalpa De Méndez'
$ws.Range("B1428").Value = 'Soto La Marina'
$ws.Range("B1437").Value = 'Amaxac De Guerrero'
$ws.Range("B1446").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B1449").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("B1452").Value = 'Papalotla De Xicohténcatl'
$ws.Range("B1453").Value = 'San Pablo Del Monte'
$ws.Range("B1458").Value = 'Tetla De La Solidaridad'
$ws.Range("B1475").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B1479").Value = 'Amatlán De Los Reyes'
$ws.Range("B1487").Value = 'Boca Del Río'
$ws.Range("B1492").Value = 'Castillo De Teayo'
$ws.Range("B1509").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1510").Value = 'Cosautlán De Carvajal'
$ws.Range("B1526").Value = 'Hueyapan De Ocampo'
$ws.Range("B1527").Value = 'Ignacio De La Llave'
$ws.Range("B1531").Value = 'Ixhuatlán De Madero'
$ws.Range("B1532").Value = 'Ixhuatlán Del Café'
$ws.Range("B1533").Value = 'Ixhuatlán Del Sureste'
$ws.Range("B1544").Value = 'Juchique De Ferrer'
$ws.Range("B1547").Value = 'Lerdo De Tejada'
$ws.Range("B1551").Value = 'Martínez De La Torre'
$ws.Range("B1553").Value = 'Medellín De Bravo'
$ws.Range("B1557").Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range("B1567").Value = 'Ozuluama De Mascareñas'
$ws.Range("B1571").Value = 'Paso De Ovejas'
$ws.Range("B1572").Value = 'Paso Del Macho'
$ws.Range("B1576").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1585").Value = 'Sayula De Alemán'
$ws.Range("B1588").Value = 'Soledad De Doblado'
$ws.Range("B1594").Value = 'Tatahuicapan De Juárez'
$ws.Range("B1614").Value = 'Tlacotepec De Mejía'
$ws.Range("B1627").Value = 'Vega De Alatorre'
$ws.Range("B1638").Value = 'Zontecomatlán De López Y Fuentes'
$ws.Range("B1639").Value = 'Zozocolco De Hidalgo'
$ws.Range("B1654").Value = 'Concepción Del Oro'
$ws.Range("B1663").Value = 'Jiménez Del Teul'
$ws.Range("B1670").Value = 'Mezquital Del Oro'
$ws.Range("B1673").Value = 'Moyahua De Estrada'
$ws.Range("B1674").Value = 'Nochistlán De Mejía'
$ws.Range("B1675").Value = 'Noria De Ángeles'
$ws.Range("B1685").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1688").Value = 'Villa De Cos'

# Minor floating point recalculation drift (last-bit rounding) on percentage column
$ws.Range("D32").Value = 0.0009764433052605884
$ws.Range("D323").Value = 0.0009764433052605884
$ws.Range("D388").Value = 0.0009764433052605884
$ws.Range("D614").Value = 0.0009764433052605884
$ws.Range("D661").Value = 0.0009764433052605884
$ws.Range("D714").Value = 0.0009764433052605884
$ws.Range("D1113").Value = 0.0009764433052605884
$ws.Range("D1130").Value = 0.09190772610765288

# Remove trailing footer/metadata rows (1697-1701); row 1696 already blank
$ws.Rows("1697:1701").Delete()

Write-Host "Edit complete"
